# This script applies the edit described by the diff:
#  - Adds 8 new data rows (new rows 5-12) with additional fruit/coffee
#    products, shifting the previous single totals row down to row 13.
#  - Updates the existing rows 2-4 content (prices / codes / brands).
#  - Recalculates the totals row (K13/L13).
#  - Extends the two color-scale conditional formats to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows before the old totals row (old row 5) so the sheet
# grows from 5 rows to 13 rows; the new rows inherit the data-row style
# (same as row 4) and the totals row is pushed down to row 13.
$ws.Range("A5:A12").EntireRow.Insert()

# Helper to force a numeric-looking string (e.g. "0253") to be stored as
# TEXT instead of being auto-converted to a number, while keeping the
# original (non quote-prefixed) cell style, by copying the number format
# from a plain text cell on the same row.
function Set-CodeText($cell, $text, $styleDonor) {
    $cell.Value = "'" + $text
    $styleDonor.Copy() | Out-Null
    $cell.PasteSpecial(-4122)
}

# --- Row 2 ---
$ws.Range("A2").Value = 'Πωλήσεις Έκπτωση 1'
$ws.Range("C2").Value = "07/16/2020"
$ws.Range("D2").Value = "07/31/2020"
$ws.Range("E2").Value = 'Nivea® Sun Spray Protect & Bronze {20} 200ml'
Set-CodeText $ws.Range("F2") '4005808859634' $ws.Range("E2")
$ws.Range("G2").Value = 8.9
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 'Nivea'
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# --- Row 3 ---
$ws.Range("A3").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C3").Value = "07/16/2020"
$ws.Range("D3").Value = "07/31/2020"
$ws.Range("E3").Value = 'Βερύκοκα® Ελληνικά (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F3") '0253' $ws.Range("E3")
$ws.Range("G3").Value = 1.65
$ws.Range("H3").Value = 1.65
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 'Βερύκοκα'
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0

# --- Row 4 ---
$ws.Range("A4").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C4").Value = "07/16/2020"
$ws.Range("D4").Value = "07/31/2020"
$ws.Range("E4").Value = 'Μπανάνες® ΕΚΟΥΑΔΟΡ (Ζυγιζόμενο) / Kgr'
Set-CodeText $ws.Range("F4") '0201' $ws.Range("E4")
$ws.Range("G4").Value = 1.25
$ws.Range("H4").Value = 1.25
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 'Μπανάνες'
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

# --- Row 5 ---
$ws.Range("A5").Value = 'Πωλήσεις Έκπτωση 1'
$ws.Range("C5").Value = "07/16/2020"
$ws.Range("D5").Value = "07/31/2020"
$ws.Range("E5").Value = 'Παπαγάλος® Ελληνικός Καφές Κουπάτος 143gr'
Set-CodeText $ws.Range("F5") '5201219486417' $ws.Range("E5")
$ws.Range("G5").Value = 3.7
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 'Παπαγάλος'
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0

# --- Row 6 ---
$ws.Range("A6").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C6").Value = "07/16/2020"
$ws.Range("D6").Value = "07/31/2020"
$ws.Range("E6").Value = 'Ροδάκινα® Ναουσας (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F6") '0208' $ws.Range("E6")
$ws.Range("G6").Value = 0.85
$ws.Range("H6").Value = 0.85
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 'Ροδάκινο'
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

# --- Row 7 ---
$ws.Range("A7").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C7").Value = "07/16/2020"
$ws.Range("D7").Value = "07/31/2020"
$ws.Range("E7").Value = 'Ροδάκινα® Ναουσας (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F7") '0208' $ws.Range("E7")
$ws.Range("G7").Value = 0.85
$ws.Range("H7").Value = 0.85
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 'Ροδάκινο'
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

# --- Row 8 ---
$ws.Range("A8").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C8").Value = "07/16/2020"
$ws.Range("D8").Value = "07/31/2020"
$ws.Range("E8").Value = 'Ροδάκινα® Ναουσας (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F8") '0208' $ws.Range("E8")
$ws.Range("G8").Value = 0.85
$ws.Range("H8").Value = 0.85
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 'Ροδάκινο'
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0

# --- Row 9 ---
$ws.Range("A9").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C9").Value = "07/16/2020"
$ws.Range("D9").Value = "07/31/2020"
$ws.Range("E9").Value = 'Βερύκοκα® Ελληνικά (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F9") '0253' $ws.Range("E9")
$ws.Range("G9").Value = 1.65
$ws.Range("H9").Value = 1.65
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 'Βερύκοκα'
$ws.Range("K9").Value = 1.7
$ws.Range("L9").Value = 2.48

# --- Row 10 ---
$ws.Range("A10").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C10").Value = "07/16/2020"
$ws.Range("D10").Value = "07/31/2020"
$ws.Range("E10").Value = 'Nutella® Πραλίνα Βάζο 400gr'
Set-CodeText $ws.Range("F10") '80135876' $ws.Range("E10")
$ws.Range("G10").Value = 3.78
$ws.Range("H10").Value = 2.95
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 'Nutella'
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 5.22

# --- Row 11 ---
$ws.Range("A11").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C11").Value = "07/16/2020"
$ws.Range("D11").Value = "07/31/2020"
$ws.Range("E11").Value = 'Μπανάνες® ΕΚΟΥΑΔΟΡ (Ζυγιζόμενο) / Kgr'
Set-CodeText $ws.Range("F11") '0201' $ws.Range("E11")
$ws.Range("G11").Value = 1.25
$ws.Range("H11").Value = 1.25
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 'Μπανάνες'
$ws.Range("K11").Value = 4.45
$ws.Range("L11").Value = 5.18

# --- Row 12 ---
$ws.Range("A12").Value = 'Πελάτες Τιμή Πώλησης'
$ws.Range("C12").Value = "07/16/2020"
$ws.Range("D12").Value = "07/31/2020"
$ws.Range("E12").Value = 'Νεκταρίνια® ΝΑΟΥΣΑΣ  (Ζυγιζόμενο) /Kgr'
Set-CodeText $ws.Range("F12") '0214' $ws.Range("E12")
$ws.Range("G12").Value = 0.85
$ws.Range("H12").Value = 0.85
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 'Νεκταρίνια'
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 5.62

# --- Totals row ---
$ws.Range("K13").Value = 14.15
$ws.Range("L13").Value = 18.5

# --- Extend conditional formatting ranges to the new data extent ---
$ws.Range("I1").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I1:I12"))
$ws.Range("J1").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J1:J12"))
